# Add placeholder values for lesson C
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$specimens = @(1, 2, 3)
$startRow = 9

for ($i = 0; $i -lt $specimens.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = "C"
    $ws.Cells.Item($row, 2).Value = $specimens[$i]
    for ($col = 3; $col -le 14; $col++) {
        $ws.Cells.Item($row, $col).Value = 0
    }
}

$ws.Range("D9").Select()
